# Generate Report for Handoff
# Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
# and the handoff/handback timestamps advance a few seconds, on all three
# sheets (Overview summary + the per-locale zh-cn / de-de detail sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 02:58:56"

# --- zh-cn detail sheet -------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 02:58:51"

# --- de-de detail sheet -------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 02:58:56"

# --- Column widths: the Status columns got narrower once the longest
# status string in use shrank from "Handed back: in sync with en-US" to
# "Ready for handoff".
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
